# Apply cell updates for the cryptos worksheet refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue($Cell, $Text) {
    # Force the value to be stored as literal text, even when it
    # looks like a number (keeps leading/trailing zeros, avoids
    # scientific notation, matches the original inlineStr cells).
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

$ws.Range("D2").Value = "62.188.80"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "2.441.07"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextCellValue $ws.Range("D5") "566.59"
$ws.Range("E5").Value = "  -1.94%  "
Set-TextCellValue $ws.Range("D6") "144.99"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("E10").Value = "  -0.12%  "
Set-TextCellValue $ws.Range("D11") "5.18"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("E12").Value = "  -3.08%  "
Set-TextCellValue $ws.Range("D13") "28.49"
$ws.Range("E13").Value = "  -2.85%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCellValue $ws.Range("D14") "0.0000172"
$ws.Range("E14").Value = "  -4.18%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.885.01"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "62.281.20"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "2.439.53"
$ws.Range("E17").Value = "  -1.57%  "
Set-TextCellValue $ws.Range("D18") "7.70"
$ws.Range("E18").Value = "  -3.21%  "
$ws.Range("E19").Value = "  -4.02%  "
Set-TextCellValue $ws.Range("D20") "319.56"
$ws.Range("E20").Value = "  -3.41%  "
Set-TextCellValue $ws.Range("D21") "4.10"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  -3.64%  "
Set-TextCellValue $ws.Range("D23") "0.999"
$ws.Range("E23").Value = "  -0.15%  "
Set-TextCellValue $ws.Range("D24") "9.80"
$ws.Range("E24").Value = "  +6.68%  "
Set-TextCellValue $ws.Range("D25") "64.83"
$ws.Range("E25").Value = "  -2.38%  "
Set-TextCellValue $ws.Range("D26") "637.14"
$ws.Range("E26").Value = "  -5.10%  "
$ws.Range("D27").Value = "2.561.06"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").Value = "0.0₃0945"
$ws.Range("E28").Value = "  -6.55%  "
Set-TextCellValue $ws.Range("D29") "0.992"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("E30").Value = "  -5.50%  "
Set-TextCellValue $ws.Range("D31") "7.79"
$ws.Range("E31").Value = "  -4.89%  "
Set-TextCellValue $ws.Range("D32") "1.80"
$ws.Range("E32").Value = "  -4.39%  "
$ws.Range("E33").Value = "  -5.47%  "
$ws.Range("E34").Value = "  -0.05%  "
Set-TextCellValue $ws.Range("D35") "1.47"
$ws.Range("E35").Value = "  -4.82%  "
$ws.Range("E36").Value = "  -4.15%  "
Set-TextCellValue $ws.Range("D37") "150.44"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("E39").Value = "  -2.93%  "
Set-TextCellValue $ws.Range("D40") "5.20"
$ws.Range("E40").Value = "  -6.30%  "
Set-TextCellValue $ws.Range("D41") "2.69"
$ws.Range("E41").Value = "  -2.62%  "
Set-TextCellValue $ws.Range("D42") "1.69"
$ws.Range("E42").Value = "  -4.87%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.74%  "
Set-TextCellValue $ws.Range("D45") "151.25"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("E46").Value = "  +0.92%  "
Set-TextCellValue $ws.Range("D47") "3.50"
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("E48").Value = "  -1.50%  "
Set-TextCellValue $ws.Range("D49") "19.92"
$ws.Range("E49").Value = "  -4.99%  "
Set-TextCellValue $ws.Range("D50") "0.0499"
$ws.Range("E50").Value = "  -3.80%  "
Set-TextCellValue $ws.Range("D51") "0.0898"
$ws.Range("E51").Value = "  -2.72%  "
